# Update the BECbIC "Employee Compensation" data row with refreshed source
# figures (file updates from rmi sep 20).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BECbIC")

$ws.Range("B2").Value = 542034080.9480869
$ws.Range("E2").Value = 63520927.93147575
$ws.Range("F2").Value = 6805554.769947206
$ws.Range("G2").Value = 1179231536.064481
$ws.Range("H2").Value = 72702948.914887
$ws.Range("I2").Value = 924819952.5592971
$ws.Range("J2").Value = 404667083.5739595
$ws.Range("L2").Value = 241695475.3087335
$ws.Range("M2").Value = 39421949.36341578
$ws.Range("N2").Value = 245202978.3234529
$ws.Range("O2").Value = 78313020.3421967
$ws.Range("P2").Value = 138121827.4900649
$ws.Range("Q2").Value = 60163263.1254696
$ws.Range("R2").Value = 81664240.49406464
$ws.Range("S2").Value = 717326256.7620312
$ws.Range("T2").Value = 2279383943.526087
$ws.Range("U2").Value = 141585829.0893514
$ws.Range("V2").Value = 893316871.9682038
$ws.Range("W2").Value = 216728079.59555
$ws.Range("X2").Value = 246596634.7808403
$ws.Range("Y2").Value = 547048827.4493796
$ws.Range("AA2").Value = 0
$ws.Range("AB2").Value = 266148772.7327911
$ws.Range("AC2").Value = 4990287427.674742
$ws.Range("AD2").Value = 4439926047.998323
$ws.Range("AF2").Value = 3223044036.876666
$ws.Range("AG2").Value = 2034365443.42734
$ws.Range("AK2").Value = 1142531647.492691
$ws.Range("AL2").Value = 11047764468.05769
$ws.Range("AM2").Value = 74442191000
$ws.Range("AN2").Value = 906098232.287804
$ws.Range("AO2").Value = 11171914285.37612
$ws.Range("AP2").Value = 671043167.4867262
$ws.Range("AQ2").Value = 0
